$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update clan overview data with latest snapshot values (serverStatusCoC en checker toegevoegd)

$ws.Range("E3").Value = 5544.0
$ws.Range("G3").Value = 9849.0

$ws.Range("G4").Value = 1392.0
$ws.Range("H4").Value = 1350.0

$ws.Range("E7").Value = 5119.0

$ws.Range("E8").Value = 5043.0

$ws.Range("G10").Value = 5451.0
$ws.Range("I10").Value = 1.98

$ws.Range("E11").Value = 4923.0

$ws.Range("E12").Value = 4846.0
$ws.Range("G12").Value = 3278.0
$ws.Range("H12").Value = 3489.0

$ws.Range("E14").Value = 4824.0
$ws.Range("H14").Value = 1661.0
$ws.Range("I14").Value = 0.69

$ws.Range("G15").Value = 1488.0
$ws.Range("H15").Value = 2085.0

$ws.Range("H24").Value = 1587.0

$ws.Range("B28").Value = "peter"
$ws.Range("C28").Value = "#8LV09JG9"
$ws.Range("D28").Value = 159.0
$ws.Range("E28").Value = 4308.0
$ws.Range("G28").Value = 2167.0
$ws.Range("H28").Value = 3146.0
$ws.Range("I28").Value = 0.69

$ws.Range("B29").Value = "elandro"
$ws.Range("C29").Value = "#22GU992L"
$ws.Range("D29").Value = 175.0
$ws.Range("E29").Value = 4285.0
$ws.Range("G29").Value = 2989.0
$ws.Range("H29").Value = 3756.0
$ws.Range("I29").Value = 0.8

$ws.Range("E31").Value = 4033.0
$ws.Range("G31").Value = 424.0
$ws.Range("H31").Value = 406.0
$ws.Range("I31").Value = 1.04

$ws.Range("E32").Value = 4030.0

$ws.Range("G35").Value = 7185.0
$ws.Range("I35").Value = 1.27

$ws.Range("G38").Value = 638.0
$ws.Range("I38").Value = 1.64

$ws.Range("E39").Value = 3424.0

$ws.Range("E41").Value = 3124.0
$ws.Range("H41").Value = 914.0
$ws.Range("I41").Value = 1.39

$ws.Range("E42").Value = 2478.0

$ws.PageSetup.LeftFooter = "Clanoverzicht"
$ws.PageSetup.RightFooter = "20/12/2017 15:01"
